# Fill in the previously-blank Time Log row 76 on Sheet1 with a new
# "Coding" time entry, mirroring what the TableView feature now writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A76: Date (10/6/2014 serial 41918)
$ws.Range("A76").Value = 41918
# B76: Start Time
$ws.Range("B76").Value = 0.93125000000000002
# C76: Stop Time
$ws.Range("C76").Value = 1.1340277777777779
# D76: Interruption (minutes)
$ws.Range("D76").Value = 30
# E76 already carries the shared formula for the Delta column; re-assert it
# so the cached result is recomputed for the newly populated row.
$ws.Range("E76").Formula = "=IF(AND(NOT(ISBLANK(B76)),NOT(ISBLANK(C76))), (C76-B76) * 24 - D76/60, """")"
# F76: Activity/category
$ws.Range("F76").Value = "Coding"

# Move the active selection to C77, matching the user's next click.
$ws.Range("C77").Select()

$excel.Calculate()
